$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1682847896440129
$ws.Range("C2").Value = 0.598705501618123
$ws.Range("J2").Value = 0.03236245954692556
$ws.Range("P2").Value = 0.1359223300970874
$ws.Range("S2").Value = 0.06472491909385113
$ws.Range("C3").Value = 0.05
$ws.Range("J3").Value = 0.04
$ws.Range("P3").Value = 0.725
$ws.Range("S3").Value = 0.185
$ws.Range("J4").Value = 0.025
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.275
$ws.Range("B6").Value = 0.05
$ws.Range("D6").Value = 0.003846153846153846
$ws.Range("F6").Value = 0.07307692307692308
$ws.Range("J6").Value = 0.2384615384615385
$ws.Range("O6").Value = 0.01538461538461539
$ws.Range("Q6").Value = 0.1692307692307692
$ws.Range("R6").Value = 0.08846153846153847
$ws.Range("S6").Value = 0.3615384615384615
$ws.Range("B7").Value = 0.0989010989010989
$ws.Range("D7").Value = 0.02197802197802198
$ws.Range("E7").Value = 0.003663003663003663
$ws.Range("F7").Value = 0.03663003663003663
$ws.Range("J7").Value = 0.1208791208791209
$ws.Range("O7").Value = 0.003663003663003663
$ws.Range("Q7").Value = 0.1758241758241758
$ws.Range("R7").Value = 0.1025641025641026
$ws.Range("S7").Value = 0.4358974358974359
$ws.Range("B8").Value = 0.1111111111111111
$ws.Range("D8").Value = 0.009578544061302681
$ws.Range("F8").Value = 0.08237547892720307
$ws.Range("J8").Value = 0.08812260536398467
$ws.Range("O8").Value = 0.01532567049808429
$ws.Range("Q8").Value = 0.1647509578544061
$ws.Range("R8").Value = 0.1168582375478927
$ws.Range("S8").Value = 0.4118773946360153
$ws.Range("B9").Value = 0.1098901098901099
$ws.Range("D9").Value = 0.02197802197802198
$ws.Range("F9").Value = 0.05494505494505494
$ws.Range("J9").Value = 0.07142857142857142
$ws.Range("O9").Value = 0.02747252747252747
$ws.Range("Q9").Value = 0.1868131868131868
$ws.Range("R9").Value = 0.1098901098901099
$ws.Range("S9").Value = 0.4175824175824176
$ws.Range("B10").Value = 0.1001494768310912
$ws.Range("D10").Value = 0.01868460388639761
$ws.Range("E10").Value = 0.0007473841554559044
$ws.Range("F10").Value = 0.07324364723467862
$ws.Range("J10").Value = 0.1270553064275037
$ws.Range("O10").Value = 0.02167414050822123
$ws.Range("Q10").Value = 0.1943198804185351
$ws.Range("R10").Value = 0.09641255605381166
$ws.Range("S10").Value = 0.3677130044843049
$ws.Range("G11").Value = 0.1465968586387434
$ws.Range("J11").Value = 0.08900523560209424
$ws.Range("K11").Value = 0.1675392670157068
$ws.Range("L11").Value = 0.5916230366492147
$ws.Range("S11").Value = 0.005235602094240838
$ws.Range("G12").Value = 0.7402597402597403
$ws.Range("J12").Value = 0.1774891774891775
$ws.Range("L12").Value = 0.03463203463203463
$ws.Range("S12").Value = 0.04761904761904762
$ws.Range("G13").Value = 0.7903225806451613
$ws.Range("J13").Value = 0.1774193548387097
$ws.Range("S13").Value = 0.03225806451612903
$ws.Range("F15").Value = 0.03555555555555556
$ws.Range("H15").Value = 0.2088888888888889
$ws.Range("I15").Value = 0.02666666666666667
$ws.Range("J15").Value = 0.3555555555555556
$ws.Range("K15").Value = 0.09777777777777778
$ws.Range("M15").Value = 0.008888888888888889
$ws.Range("O15").Value = 0.04444444444444445
$ws.Range("S15").Value = 0.2222222222222222
$ws.Range("F16").Value = 0.01428571428571429
$ws.Range("H16").Value = 0.119047619047619
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.4619047619047619
$ws.Range("K16").Value = 0.1571428571428571
$ws.Range("M16").Value = 0.01428571428571429
$ws.Range("N16").Value = 0.004761904761904762
$ws.Range("O16").Value = 0.03333333333333333
$ws.Range("S16").Value = 0.1238095238095238
$ws.Range("F17").Value = 0.0128755364806867
$ws.Range("H17").Value = 0.1931330472103004
$ws.Range("I17").Value = 0.07939914163090128
$ws.Range("J17").Value = 0.4356223175965665
$ws.Range("K17").Value = 0.1072961373390558
$ws.Range("M17").Value = 0.02789699570815451
$ws.Range("O17").Value = 0.06008583690987124
$ws.Range("S17").Value = 0.08369098712446352
$ws.Range("F18").Value = 0.015625
$ws.Range("H18").Value = 0.1875
$ws.Range("I18").Value = 0.08203125
$ws.Range("J18").Value = 0.44140625
$ws.Range("K18").Value = 0.109375
$ws.Range("M18").Value = 0.02734375
$ws.Range("N18").Value = 0.00390625
$ws.Range("O18").Value = 0.06640625
$ws.Range("S18").Value = 0.06640625
$ws.Range("F19").Value = 0.01878287002253944
$ws.Range("H19").Value = 0.2344102178812923
$ws.Range("I19").Value = 0.0811419984973704
$ws.Range("J19").Value = 0.33809166040571
$ws.Range("K19").Value = 0.1292261457550714
$ws.Range("M19").Value = 0.02930127723516153
$ws.Range("N19").Value = 0.001502629601803156
$ws.Range("O19").Value = 0.067618332081142
$ws.Range("S19").Value = 0.09992486851990984
